$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed cells keep their original "text" storage (they currently
# hold numeric-looking strings such as "306.24" or "6.33%"). Pre-marking each
# cell as Text before writing the new value prevents Excel from re-interpreting
# the new numeric-looking string as a Number/Percentage.
$cellsToUpdate = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($addr in $cellsToUpdate) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "307.25"
$ws.Range("E2").Value = "6.81%"
$ws.Range("D3").Value = "35.12"
$ws.Range("E3").Value = "13.16%"
$ws.Range("D4").Value = "5.171"
$ws.Range("E4").Value = "5.26%"
$ws.Range("D5").Value = "0.07904"
$ws.Range("E5").Value = "7.92%"
$ws.Range("D6").Value = "2.359"
$ws.Range("E6").Value = "5.87%"
$ws.Range("D7").Value = "8.018"
$ws.Range("E7").Value = "3.61%"
$ws.Range("D8").Value = "3.975"
$ws.Range("E8").Value = "6.82%"
$ws.Range("D9").Value = "0.9287"
$ws.Range("E9").Value = "2.73%"
$ws.Range("E10").Value = "11.43%"
$ws.Range("D11").Value = "0.1836"
$ws.Range("E11").Value = "8.37%"
$ws.Range("D12").Value = "0.08668"
$ws.Range("E12").Value = "4.73%"
$ws.Range("D13").Value = "0.03400"
$ws.Range("E13").Value = "8.77%"
$ws.Range("D14").Value = "0.09881"
$ws.Range("E14").Value = "-0.47%"
$ws.Range("D15").Value = "0.001484"
$ws.Range("E15").Value = "-1.28%"
$ws.Range("D16").Value = "0.005665"
$ws.Range("E16").Value = "-0.77%"
$ws.Range("D17").Value = "3.502"
$ws.Range("E17").Value = "-0.94%"
$ws.Range("D18").Value = "2.110"
$ws.Range("E18").Value = "1.35%"
$ws.Range("D19").Value = "0.3399"
$ws.Range("E19").Value = "2.06%"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").Value = "0.29%"
$ws.Range("D21").Value = "4.549"
$ws.Range("E21").Value = "8.43%"
$ws.Range("E22").Value = "8.98%"
$ws.Range("D23").Value = "0.04563"
$ws.Range("E23").Value = "1.22%"
$ws.Range("D24").Value = "0.001214"
$ws.Range("E24").Value = "0.33%"
$ws.Range("D25").Value = "0.004479"
$ws.Range("E25").Value = "7.83%"
$ws.Range("D26").Value = "0.0001291"
$ws.Range("E26").Value = "-0.76%"
$ws.Range("D27").Value = "0.0003388"
$ws.Range("E27").Value = "-0.27%"
$ws.Range("D39").Value = "0.01795"
$ws.Range("E39").Value = "14.22%"
$ws.Range("D40").Value = "0.04804"
$ws.Range("E40").Value = "8.06%"
$ws.Range("D41").Value = "0.007794"
$ws.Range("E41").Value = "6.25%"
$ws.Range("D42").Value = "0.1425"
$ws.Range("E42").Value = "7.39%"
$ws.Range("D43").Value = "0.007056"
$ws.Range("E43").Value = "-26.11%"
$ws.Range("D44").Value = "0.002202"
$ws.Range("E44").Value = "-3.90%"
$ws.Range("D45").Value = "0.009602"
$ws.Range("E45").Value = "15.38%"
$ws.Range("D46").Value = "0.00005980"
$ws.Range("E46").Value = "-2.26%"
$ws.Range("D47").Value = "0.00000000745"
$ws.Range("E47").Value = "-0.76%"
$ws.Range("E48").Value = "13.25%"
$ws.Range("D49").Value = "0.002680"
$ws.Range("E49").Value = "33.87%"
$ws.Range("D50").Value = "0.00002086"
$ws.Range("E50").Value = "-0.76%"
$ws.Range("D51").Value = "0.0001986"
$ws.Range("E51").Value = "-0.76%"

Write-Host "Updated symbol list"
